$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gUSD 26.06.25")
$ws.Activate()

$ws.Cells.Item(150, 3).Value = 25.3655
$ws.Cells.Item(150, 4).Value = 102.668
$ws.Cells.Item(150, 5).Value = 9.42
$ws.Cells.Item(150, 6).Value = 7.27
$ws.Cells.Item(150, 7).Value = 10.06
$ws.Cells.Item(150, 9).Formula = "=C150/`$D`$3"
$ws.Cells.Item(150, 9).NumberFormat = "0.0000"
$ws.Cells.Item(150, 13).Formula = "=C150+D150"

$ws.Cells.Item(151, 3).Value = 24.0189
$ws.Cells.Item(151, 4).Value = 103.826
$ws.Cells.Item(151, 5).Value = 9.38
$ws.Cells.Item(151, 6).Value = 6.53
$ws.Cells.Item(151, 7).Value = 8.85
$ws.Cells.Item(151, 9).Formula = "=C151/`$D`$3"
$ws.Cells.Item(151, 9).NumberFormat = "0.0000"
$ws.Cells.Item(151, 13).Formula = "=C151+D151"

$ws.Cells.Item(152, 3).Value = 22.6624
$ws.Cells.Item(152, 4).Value = 104.268
$ws.Cells.Item(152, 5).Value = 9.3
$ws.Cells.Item(152, 6).Value = 6.12
$ws.Cells.Item(152, 7).Value = 3.56
$ws.Cells.Item(152, 9).Formula = "=C152/`$D`$3"
$ws.Cells.Item(152, 9).NumberFormat = "0.0000"
$ws.Cells.Item(152, 13).Formula = "=C152+D152"

$ws.Cells.Item(153, 3).Value = 19.9926
$ws.Cells.Item(153, 4).Value = 104.784
$ws.Cells.Item(153, 5).Value = 8.71
$ws.Cells.Item(153, 6).Value = 6.17
$ws.Cells.Item(153, 7).Value = 3.45
$ws.Cells.Item(153, 9).Formula = "=C153/`$D`$3"
$ws.Cells.Item(153, 9).NumberFormat = "0.0000"
$ws.Cells.Item(153, 13).Formula = "=C153+D153"

$ws.Cells.Item(154, 3).Value = 19.1447
$ws.Cells.Item(154, 4).Value = 105.801
$ws.Cells.Item(154, 5).Value = 8.86
$ws.Cells.Item(154, 6).Value = 6.69
$ws.Cells.Item(154, 7).Value = 7.8
$ws.Cells.Item(154, 9).Formula = "=C154/`$D`$3"
$ws.Cells.Item(154, 9).NumberFormat = "0.0000"
$ws.Cells.Item(154, 13).Formula = "=C154+D154"

$ws.Cells.Item(155, 3).Value = 15.881
$ws.Cells.Item(155, 4).Value = 107.202
$ws.Cells.Item(155, 5).Value = 7.83
$ws.Cells.Item(155, 6).Value = 7
$ws.Cells.Item(155, 7).Value = 10.21
$ws.Cells.Item(155, 9).Formula = "=C155/`$D`$3"
$ws.Cells.Item(155, 9).NumberFormat = "0.0000"
$ws.Cells.Item(155, 13).Formula = "=C155+D155"

$ws.Cells.Item(156, 3).Value = 18.2183
$ws.Cells.Item(156, 4).Value = 108.979
$ws.Cells.Item(156, 5).Value = 9.7
$ws.Cells.Item(156, 6).Value = 8.25
$ws.Cells.Item(156, 7).Value = 14.03
$ws.Cells.Item(156, 9).Formula = "=C156/`$D`$3"
$ws.Cells.Item(156, 9).NumberFormat = "0.0000"
$ws.Cells.Item(156, 13).Formula = "=C156+D156"

$ws.Cells.Item(157, 3).Value = 16.694
$ws.Cells.Item(157, 4).Value = 111.456
$ws.Cells.Item(157, 5).Value = 9.66
$ws.Cells.Item(157, 6).Value = 9.47
$ws.Cells.Item(157, 7).Value = 19.04
$ws.Cells.Item(157, 9).Formula = "=C157/`$D`$3"
$ws.Cells.Item(157, 9).NumberFormat = "0.0000"
$ws.Cells.Item(157, 13).Formula = "=C157+D157"

$ws.Cells.Item(158, 3).Value = 15.5038
$ws.Cells.Item(158, 4).Value = 114.613
$ws.Cells.Item(158, 5).Value = 9.68
$ws.Cells.Item(158, 6).Value = 11.82
$ws.Cells.Item(158, 7).Value = 26.48
$ws.Cells.Item(158, 9).Formula = "=C158/`$D`$3"
$ws.Cells.Item(158, 9).NumberFormat = "0.0000"
$ws.Cells.Item(158, 13).Formula = "=C158+D158"

$ws.Cells.Item(159, 3).Value = 13.8024
$ws.Cells.Item(159, 4).Value = 115.752
$ws.Cells.Item(159, 5).Value = 9.39
$ws.Cells.Item(159, 6).Value = 12.56
$ws.Cells.Item(159, 7).Value = 8.9
$ws.Cells.Item(159, 9).Formula = "=C159/`$D`$3"
$ws.Cells.Item(159, 9).NumberFormat = "0.0000"
$ws.Cells.Item(159, 13).Formula = "=C159+D159"

$ws.Cells.Item(160, 3).Value = 12.7746
$ws.Cells.Item(160, 4).Value = 116.276
$ws.Cells.Item(160, 5).Value = 9.51
$ws.Cells.Item(160, 6).Value = 12.66
$ws.Cells.Item(160, 7).Value = 4.18
$ws.Cells.Item(160, 9).Formula = "=C160/`$D`$3"
$ws.Cells.Item(160, 9).NumberFormat = "0.0000"
$ws.Cells.Item(160, 13).Formula = "=C160+D160"

$ws.Cells.Item(161, 3).Value = 11.6147
$ws.Cells.Item(161, 4).Value = 117.459
$ws.Cells.Item(161, 5).Value = 9.63
$ws.Cells.Item(161, 6).Value = 12.85
$ws.Cells.Item(161, 7).Value = 8.94
$ws.Cells.Item(161, 9).Formula = "=C161/`$D`$3"
$ws.Cells.Item(161, 9).NumberFormat = "0.0000"
$ws.Cells.Item(161, 13).Formula = "=C161+D161"

$ws.Cells.Item(162, 3).Value = 10.4317
$ws.Cells.Item(162, 4).Value = 118.444
$ws.Cells.Item(162, 5).Value = 9.74
$ws.Cells.Item(162, 6).Value = 12.41
$ws.Cells.Item(162, 7).Value = 7.39
$ws.Cells.Item(162, 9).Formula = "=C162/`$D`$3"
$ws.Cells.Item(162, 9).NumberFormat = "0.0000"
$ws.Cells.Item(162, 13).Formula = "=C162+D162"

$ws.Cells.Item(163, 3).Value = 9.29227
$ws.Cells.Item(163, 4).Value = 119.597
$ws.Cells.Item(163, 5).Value = 9.94
$ws.Cells.Item(163, 6).Value = 11.65
$ws.Cells.Item(163, 7).Value = 8.66
$ws.Cells.Item(163, 9).Formula = "=C163/`$D`$3"
$ws.Cells.Item(163, 9).NumberFormat = "0.0000"
$ws.Cells.Item(163, 13).Formula = "=C163+D163"

$ws.Cells.Item(164, 3).Value = 8.18784
$ws.Cells.Item(164, 4).Value = 120.164
$ws.Cells.Item(164, 5).Value = 10.18
$ws.Cells.Item(164, 6).Value = 10.02
$ws.Cells.Item(164, 7).Value = 4.9
$ws.Cells.Item(164, 9).Formula = "=C164/`$D`$3"
$ws.Cells.Item(164, 9).NumberFormat = "0.0000"
$ws.Cells.Item(164, 13).Formula = "=C164+D164"

$ws.Cells.Item(165, 3).Value = 7.04746
$ws.Cells.Item(165, 4).Value = 123.389
$ws.Cells.Item(165, 5).Value = 10.57
$ws.Cells.Item(165, 6).Value = 9.67
$ws.Cells.Item(165, 7).Value = 26.34
$ws.Cells.Item(165, 9).Formula = "=C165/`$D`$3"
$ws.Cells.Item(165, 9).NumberFormat = "0.0000"
$ws.Cells.Item(165, 13).Formula = "=C165+D165"

$ws.Cells.Item(166, 3).Value = 5.28939
$ws.Cells.Item(166, 4).Value = 124.284
$ws.Cells.Item(166, 5).Value = 9.95
$ws.Cells.Item(166, 6).Value = 9.31
$ws.Cells.Item(166, 7).Value = 6.65
$ws.Cells.Item(166, 9).Formula = "=C166/`$D`$3"
$ws.Cells.Item(166, 9).NumberFormat = "0.0000"
$ws.Cells.Item(166, 13).Formula = "=C166+D166"

$ws.Cells.Item(167, 3).Value = 4.67842
$ws.Cells.Item(167, 4).Value = 126.243
$ws.Cells.Item(167, 5).Value = 11.93
$ws.Cells.Item(167, 6).Value = 10.87
$ws.Cells.Item(167, 7).Value = 15.08
$ws.Cells.Item(167, 9).Formula = "=C167/`$D`$3"
$ws.Cells.Item(167, 9).NumberFormat = "0.0000"
$ws.Cells.Item(167, 13).Formula = "=C167+D167"

$ws.Cells.Item(168, 3).Value = 2.12643
$ws.Cells.Item(168, 4).Value = 127.37
$ws.Cells.Item(168, 5).Value = 11.2
$ws.Cells.Item(168, 6).Value = 10.47
$ws.Cells.Item(168, 7).Value = 5.9
$ws.Cells.Item(168, 9).Formula = "=C168/`$D`$3"
$ws.Cells.Item(168, 9).NumberFormat = "0.0000"
$ws.Cells.Item(168, 13).Formula = "=C168+D168"

# Update the view: scroll position and selection to match the edited area
$excel.ActiveWindow.ScrollRow = 156
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C169").Select()

Write-Output "done"